# Update 1 ao 15
# "Projeto Hair Tech" -> split into "Projeto " / "Hair" / " Tech" so that
# "Hair" can carry its own (spell-check) run boundary, matching how Word
# itself fragments a run when it flags an embedded foreign word.
#
# All three resulting runs must keep the exact same run formatting
# (Arial, bCs, color 595959/text1/A6 tint, sz/szCs 32) that the original
# single run had - only the text is being split, nothing about the look
# changes.

$d = $word.ActiveDocument

# Locate the exact phrase and split it into three runs of identical
# formatting by toggling a boolean character property (Bold) on/off over
# just the "Hair" sub-range: turning it on forces Word to fragment the
# run at the sub-range boundaries, and turning it back off removes the
# (now redundant) explicit <w:b/> again, leaving three runs that share
# the original formatting, with "Projeto " / "Hair" / " Tech" as their
# text content.
$rng = $d.Content
$rng.Find.Execute("Hair")
$rng.Bold = 1
$rng.Bold = 0
